$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09788908052730215
$ws.Range("C2").Value = 0.335102348664815
$ws.Range("D2").Value = 0.3535955084851918
$ws.Range("E2").Value = 0.594638973230978
$ws.Range("F2").Value = 0.5922487602141151

$ws.Range("B3").Value = 0.2490523672941219
$ws.Range("C3").Value = 0.7190685700340571
$ws.Range("D3").Value = 1.215519356195628
$ws.Range("E3").Value = 1.10250594383687
$ws.Range("F3").Value = 1.084694482345899

$ws.Range("B4").Value = 0.5894233596642614
$ws.Range("C4").Value = 0.9008946303567875
$ws.Range("D4").Value = 2.205314398756898
$ws.Range("E4").Value = 1.485030100286488
$ws.Range("F4").Value = 1.376884423676235

$ws.Range("B5").Value = 0.3608719394995029
$ws.Range("C5").Value = 0.9602851306804567
$ws.Range("D5").Value = 2.593079633011949
$ws.Range("E5").Value = 1.610304205115278
$ws.Range("F5").Value = 1.585610441980187
$ws.Range("G5").Value = 49

$ws.Range("B6").Value = 0.4922773315969051
$ws.Range("C6").Value = 0.9502479493674117
$ws.Range("D6").Value = 2.535544773781611
$ws.Range("E6").Value = 1.592339402822655
$ws.Range("F6").Value = 1.53035925780564
$ws.Range("G6").Value = 48

$ws.Range("B7").Value = 0.3820685963133433
$ws.Range("C7").Value = 0.8448709767778345
$ws.Range("D7").Value = 2.538317255165756
$ws.Range("E7").Value = 1.593209733577396
$ws.Range("F7").Value = 1.566938792948644
$ws.Range("G7").Value = 39

$ws.Range("B8").Value = 0.4051237794541663
$ws.Range("C8").Value = 0.9263031120077183
$ws.Range("D8").Value = 2.696411648248952
$ws.Range("E8").Value = 1.642075408819264
$ws.Range("F8").Value = 1.612676825583574
$ws.Range("G8").Value = 38

$ws.Range("B9").Value = 0.3398840034424069
$ws.Range("C9").Value = 1.065209915832154
$ws.Range("D9").Value = 4.137733237635551
$ws.Range("E9").Value = 2.034141892208002
$ws.Range("F9").Value = 2.05507243349997
$ws.Range("G9").Value = 21

$ws.Range("B10").Value = 0.05151455505799022
$ws.Range("C10").Value = 0.8135860526754219
$ws.Range("D10").Value = 1.803889379070564
$ws.Range("E10").Value = 1.343089490343277
$ws.Range("F10").Value = 1.392764235822703
$ws.Range("G10").Value = 14

$ws.Range("B11").Value = 0.8476851521374721
$ws.Range("C11").Value = 0.8476851521374721
$ws.Range("D11").Value = 0.9146264493203358
$ws.Range("E11").Value = 0.9563610454845679
$ws.Range("F11").Value = 0.4950458718214991
